# Auto-generated Excel COM-interop script to apply numeric updates
# to columns H-N (profit calc columns) across multiple sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2899.2
$ws.Range("I32").Value = 2998
$ws.Range("J32").Value = 2833.3333
$ws.Range("K32").Value = 2998
$ws.Range("L32").Value = 2833.3333
$ws.Range("M32").Value = -2672
$ws.Range("N32").Value = -3485.3333
$ws.Range("H112").Value = 2626
$ws.Range("J112").Value = 2682.9092
$ws.Range("L112").Value = 8048.7276
$ws.Range("N112").Value = -10264.7276
$ws.Range("H132").Value = 2203.0476
$ws.Range("I132").Value = 2163.2
$ws.Range("K132").Value = 6489.599999999999
$ws.Range("M132").Value = -3959.599999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4945.375
$ws.Range("I63").Value = 2401.2856
$ws.Range("K63").Value = 2401.2856
$ws.Range("M63").Value = -1715.2856
$ws.Range("H66").Value = 4945.375
$ws.Range("I66").Value = 2401.2856
$ws.Range("K66").Value = 12006.428
$ws.Range("M66").Value = -8574.428
$ws.Range("H132").Value = 1549.3684
$ws.Range("I132").Value = 1574.3334
$ws.Range("K132").Value = 4723.0002
$ws.Range("M132").Value = -2193.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3603.55
$ws.Range("I99").Value = 3614.4614
$ws.Range("K99").Value = 3614.4614
$ws.Range("M99").Value = -2116.4614
$ws.Range("H105").Value = 5173.5884
$ws.Range("I105").Value = 3622
$ws.Range("K105").Value = 3622
$ws.Range("M105").Value = -1875
$ws.Range("H134").Value = 3163.5
$ws.Range("I134").Value = 2762
$ws.Range("K134").Value = 8286
$ws.Range("M134").Value = -5751

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2892.122
$ws.Range("I31").Value = 1503.6333
$ws.Range("K31").Value = 1503.6333
$ws.Range("M31").Value = -1208.6333
$ws.Range("H34").Value = 2892.122
$ws.Range("I34").Value = 1503.6333
$ws.Range("K34").Value = 1503.6333
$ws.Range("M34").Value = -1301.6333
$ws.Range("H99").Value = 13674.774
$ws.Range("J99").Value = 15338.315
$ws.Range("L99").Value = 15338.315
$ws.Range("N99").Value = -18334.315
$ws.Range("H122").Value = 3503.8235
$ws.Range("I122").Value = 3734
$ws.Range("K122").Value = 11202
$ws.Range("M122").Value = -8752
$ws.Range("H126").Value = 13674.774
$ws.Range("J126").Value = 15338.315
$ws.Range("L126").Value = 46014.945
$ws.Range("N126").Value = -50954.945

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1000.3333
$ws.Range("I12").Value = 999
$ws.Range("J12").Value = 1001
$ws.Range("K12").Value = 2997
$ws.Range("L12").Value = 3003
$ws.Range("M12").Value = -2824
$ws.Range("N12").Value = -3349
$ws.Range("H17").Value = 130
$ws.Range("I17").Value = 138.25
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = 414.75
$ws.Range("L17").Value = 291
$ws.Range("M17").Value = -245.75
$ws.Range("N17").Value = -629
$ws.Range("H63").Value = 137.33333
$ws.Range("I63").Value = 137.33333
$ws.Range("K63").Value = 411.99999
$ws.Range("M63").Value = 337.00001
$ws.Range("H66").Value = 137.33333
$ws.Range("I66").Value = 137.33333
$ws.Range("K66").Value = 1235.99997
$ws.Range("M66").Value = 2508.00003
$ws.Range("H69").Value = 2119.3
$ws.Range("J69").Value = 1986.625
$ws.Range("L69").Value = 5959.875
$ws.Range("N69").Value = -7581.875
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 12000
$ws.Range("M70").Value = -11685
$ws.Range("H72").Value = 2119.3
$ws.Range("J72").Value = 1986.625
$ws.Range("L72").Value = 17879.625
$ws.Range("N72").Value = -25991.625
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 12000
$ws.Range("M73").Value = -10908
$ws.Range("H117").Value = 1417.6923
$ws.Range("I117").Value = 273.33334
$ws.Range("K117").Value = 820.0000200000001
$ws.Range("M117").Value = 2621.99998
$ws.Range("H129").Value = 1310.625
$ws.Range("I129").Value = 298.2
$ws.Range("J129").Value = 2998
$ws.Range("K129").Value = 894.5999999999999
$ws.Range("L129").Value = 8994
$ws.Range("M129").Value = 4105.4
$ws.Range("N129").Value = -18994
$ws.Range("H140").Value = 2978.1667
$ws.Range("I140").Value = 2978.1667
$ws.Range("K140").Value = 8934.500100000001
$ws.Range("M140").Value = -3754.500100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I80").Value = 4166.8335
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4166.8335
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3168.8335
$ws.Range("N80").Value = ""
$ws.Range("I83").Value = 4166.8335
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 20834.1675
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -15842.1675
$ws.Range("N83").Value = ""
$ws.Range("H132").Value = 1846.1111
$ws.Range("I132").Value = 1035.2307
$ws.Range("J132").Value = 3954.4
$ws.Range("K132").Value = 3105.6921
$ws.Range("L132").Value = 11863.2
$ws.Range("M132").Value = -575.6921000000002
$ws.Range("N132").Value = -16923.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1370.9166
$ws.Range("I16").Value = 1465.7273
$ws.Range("K16").Value = 1465.7273
$ws.Range("M16").Value = -1295.7273
$ws.Range("H18").Value = 69999
$ws.Range("J18").Value = 69999
$ws.Range("L18").Value = 69999
$ws.Range("N18").Value = -70343
$ws.Range("H40").Value = 3552.5
$ws.Range("I40").Value = 3400
$ws.Range("J40").Value = 3705
$ws.Range("K40").Value = 3400
$ws.Range("L40").Value = 3705
$ws.Range("M40").Value = -3264
$ws.Range("N40").Value = -3977
$ws.Range("H46").Value = 2298.9412
$ws.Range("I46").Value = 1711.8572
$ws.Range("J46").Value = 2709.9
$ws.Range("K46").Value = 1711.8572
$ws.Range("L46").Value = 2709.9
$ws.Range("M46").Value = -1523.8572
$ws.Range("N46").Value = -3085.9
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -20980
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20466
$ws.Range("H82").Value = 3388.6956
$ws.Range("I82").Value = 3579.1177
$ws.Range("K82").Value = 3579.1177
$ws.Range("M82").Value = -3218.1177
$ws.Range("H85").Value = 3388.6956
$ws.Range("I85").Value = 3579.1177
$ws.Range("K85").Value = 3579.1177
$ws.Range("M85").Value = -2331.1177

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7248.3125
$ws.Range("I62").Value = 5897.2
$ws.Range("J62").Value = 7862.4546
$ws.Range("K62").Value = 5897.2
$ws.Range("L62").Value = 7862.4546
$ws.Range("M62").Value = -5273.2
$ws.Range("N62").Value = -9110.454600000001
$ws.Range("H65").Value = 7248.3125
$ws.Range("I65").Value = 5897.2
$ws.Range("J65").Value = 7862.4546
$ws.Range("K65").Value = 29486
$ws.Range("L65").Value = 39312.273
$ws.Range("M65").Value = -26366
$ws.Range("N65").Value = -45552.273
$ws.Range("H132").Value = 1950.5
$ws.Range("I132").Value = 1763.125
$ws.Range("K132").Value = 5289.375
$ws.Range("M132").Value = -2759.375
